$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixed bug in processing data from database
$ws.Range("E1").Value = "coin"
$ws.Range("E2").Value = 30

# Update selection to E3
$ws.Range("E3").Select()
